$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Update the "Metrics" sheet values (B2:B13) with the new figures.
# Downstream formulas on "today" (and anywhere else) that reference these
# cells will recalculate automatically.
# ---------------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value = 112073.83000000002
$metrics.Range("B3").Value = 85184.72
$metrics.Range("B4").Value = 27939.27
$metrics.Range("B5").Value = 4523
$metrics.Range("B6").Value = 5747944.5599999996
$metrics.Range("B7").Value = 4855902.3500000006
$metrics.Range("B8").Value = 1692031.09
$metrics.Range("B9").Value = 224800
$metrics.Range("B10").Value = 34213325.549999997
$metrics.Range("B11").Value = 32131177.510000002
$metrics.Range("B12").Value = 11973753.129999999
$metrics.Range("B13").Value = 1322430

# ---------------------------------------------------------------------------
# Clear the manual override formulas in B3:B6 on the "today" sheet so they
# go back to blank cells (keeping their existing style).
# ---------------------------------------------------------------------------
$today = $wb.Worksheets.Item("today")
$today.Range("B3").ClearContents()
$today.Range("B4").ClearContents()
$today.Range("B5").ClearContents()
$today.Range("B6").ClearContents()

# ---------------------------------------------------------------------------
# Switch the active/selected sheet from "Metrics" to "today" and move the
# selection there to F9 (was F11:F22).
# ---------------------------------------------------------------------------
$today.Activate()
$today.Range("F9").Select()
